$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $cell = $ws.Range($cellRef)
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = "Normal"
}

Set-TextValue 'D2' '79.593.19'
Set-TextValue 'E2' '  +4.50%  '
Set-TextValue 'D3' '3.176.70'
Set-TextValue 'E3' '  +3.52%  '
Set-TextValue 'E4' '  +0.41%  '
Set-TextValue 'D5' '207.01'
Set-TextValue 'E5' '  +4.55%  '
Set-TextValue 'D6' '625.55'
Set-TextValue 'E6' '  +1.67%  '
Set-TextValue 'D7' '0.267'
Set-TextValue 'E7' '  +28.23%  '
Set-TextValue 'D8' '1.00'
Set-TextValue 'E8' '  +0.27%  '
Set-TextValue 'D9' '0.595'
Set-TextValue 'E9' '  +8.25%  '
Set-TextValue 'D10' '3.180.14'
Set-TextValue 'E10' '  +3.71%  '
Set-TextValue 'D11' '0.611'
Set-TextValue 'E11' '  +38.95%  '
Set-TextValue 'D12' '0.0000255'
Set-TextValue 'E12' '  +32.11%  '
Set-TextValue 'D13' '0.165'
Set-TextValue 'E13' '  +2.59%  '
Set-TextValue 'B14' 'WrappedliquidstakedEther2.0'
Set-TextValue 'C14' 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
Set-TextValue 'D14' '3.798.80'
Set-TextValue 'E14' '  +5.24%  '
Set-TextValue 'B15' 'Toncoin'
Set-TextValue 'C15' 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
Set-TextValue 'D15' '5.28'
Set-TextValue 'E15' '  +1.14%  '
Set-TextValue 'D16' '31.88'
Set-TextValue 'E16' '  +10.30%  '
Set-TextValue 'D17' '80.017.11'
Set-TextValue 'E17' '  +4.97%  '
Set-TextValue 'D18' '3.217.03'
Set-TextValue 'E18' '  +5.09%  '
Set-TextValue 'D19' '14.36'
Set-TextValue 'E19' '  +5.99%  '
Set-TextValue 'D20' '440.13'
Set-TextValue 'E20' '  +16.05%  '
Set-TextValue 'D21' '9.16'
Set-TextValue 'E21' '  +0.62%  '
Set-TextValue 'D22' '2.93'
Set-TextValue 'E22' '  +19.50%  '
Set-TextValue 'D23' '5.29'
Set-TextValue 'E23' '  +20.47%  '
Set-TextValue 'D24' '3.381.47'
Set-TextValue 'E24' '  +4.65%  '
Set-TextValue 'D25' '76.83'
Set-TextValue 'E25' '  +6.66%  '
Set-TextValue 'D26' '4.70'
Set-TextValue 'E26' '  +8.45%  '
Set-TextValue 'D27' '10.84'
Set-TextValue 'E27' '  +10.12%  '
Set-TextValue 'D28' '0.999'
Set-TextValue 'E28' '  -0.13%  '
Set-TextValue 'D29' '0.0000122'
Set-TextValue 'E29' '  +13.21%  '
Set-TextValue 'D30' '9.08'
Set-TextValue 'E30' '  +9.51%  '
Set-TextValue 'D31' '1.00'
Set-TextValue 'E31' '  +0.45%  '
Set-TextValue 'B32' 'Bittensor'
Set-TextValue 'C32' 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
Set-TextValue 'D32' '544.28'
Set-TextValue 'E32' '  +9.29%  '
Set-TextValue 'B33' 'Fetch.AI'
Set-TextValue 'C33' 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
Set-TextValue 'D33' '1.48'
Set-TextValue 'E33' '  +4.20%  '
Set-TextValue 'B34' 'PancakeSwap'
Set-TextValue 'C34' 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
Set-TextValue 'D34' '2.00'
Set-TextValue 'E34' '  +4.57%  '
Set-TextValue 'B35' 'Kaspa'
Set-TextValue 'C35' 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
Set-TextValue 'D35' '0.149'
Set-TextValue 'E35' '  +20.94%  '
Set-TextValue 'D36' '23.21'
Set-TextValue 'E36' '  +12.20%  '
Set-TextValue 'D37' '0.123'
Set-TextValue 'E37' '  +20.60%  '
Set-TextValue 'B38' 'FirstDigitalUSD'
Set-TextValue 'C38' 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
Set-TextValue 'D38' '1.00'
Set-TextValue 'E38' '  +0.36%  '
Set-TextValue 'B39' 'PolygonEcosystemToken'
Set-TextValue 'C39' 'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol'
Set-TextValue 'D39' '0.408'
Set-TextValue 'E39' '  +8.37%  '
Set-TextValue 'B40' 'WhiteBITCoin'
Set-TextValue 'C40' 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
Set-TextValue 'D40' '20.77'
Set-TextValue 'E40' '  +3.64%  '
Set-TextValue 'B41' 'Monero'
Set-TextValue 'C41' 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
Set-TextValue 'D41' '164.64'
Set-TextValue 'E41' '  +1.07%  '
Set-TextValue 'D42' '5.64'
Set-TextValue 'E42' '  +10.54%  '
Set-TextValue 'D44' '188.72'
Set-TextValue 'E44' '  -2.89%  '
Set-TextValue 'D45' '1.81'
Set-TextValue 'E45' '  +10.30%  '
Set-TextValue 'D46' '2.68'
Set-TextValue 'E46' '  +10.66%  '
Set-TextValue 'D47' '0.784'
Set-TextValue 'E47' '  -1.39%  '
Set-TextValue 'B48' 'OKB'
Set-TextValue 'C48' 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
Set-TextValue 'D48' '43.64'
Set-TextValue 'E48' '  +5.55%  '
Set-TextValue 'B49' 'ImmutableX'
Set-TextValue 'C49' 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
Set-TextValue 'D49' '1.31'
Set-TextValue 'E49' '  +5.10%  '
Set-TextValue 'D50' '4.26'
Set-TextValue 'E50' '  +10.00%  '
Set-TextValue 'D51' '0.630'
Set-TextValue 'E51' '  +6.20%  '

Write-Output "Applied 121 cell updates"